# CU 25 - Consultar gastos
# Updates the "Lista de Tareas" sheet: mark the row-13/row-14 tasks as
# "Hecho" (done), fill in the estimated hours and day-3 consumption, and
# move the frozen-pane / selection to reflect the new working cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# --- Row 13 -----------------------------------------------------------
$ws.Range("F13").Value = "Hecho"
$ws.Range("G13").Value = 2
$ws.Range("N13").Value = 3

# --- Row 14 -----------------------------------------------------------
$ws.Range("F14").Value = "Hecho"
$ws.Range("G14").Value = 1
$ws.Range("N14").Value = 1

# --- View state ---------------------------------------------------------
$ws.Range("N18").Select()
$excel.ActiveWindow.FreezePanes = $true
